# Generate Report for Handoff
# Updates the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
# timestamps for the most recently processed file
# (aaa3bc37-271c-40d3-94bf-a7de80abcf60.md) across the Overview, zh-cn and
# de-de sheets, reflecting a fresh handoff xliff generation run.

$wb = $excel.ActiveWorkbook

# Overview sheet, row 7 ("aaa3bc37-271c-40d3-94bf-a7de80abcf60.md"),
# column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G7").Value = "2016-08-17 20:43:36"

# zh-cn sheet, row 7 ("aaa3bc37-271c-40d3-94bf-a7de80abcf60.md"),
# column H = "Latest Handoff Datetime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H7").Value = "2016-08-17 20:43:31"

# de-de sheet, row 7 ("aaa3bc37-271c-40d3-94bf-a7de80abcf60.md"),
# column H = "Latest Handoff Datetime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H7").Value = "2016-08-17 20:43:36"
